$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column K (최종점수) values
$ws.Range("K2").Value = 66
$ws.Range("K3").Value = 54
$ws.Range("K4").Value = 51.2
$ws.Range("K5").Value = 50
$ws.Range("K6").Value = 50

# Update column N (MACRO_SCORE) values
$ws.Range("N2").Value = 85.96878041621773
$ws.Range("N3").Value = 85.96878041621773
$ws.Range("N4").Value = 85.96878041621773
$ws.Range("N5").Value = 85.96878041621773
$ws.Range("N6").Value = 85.96878041621773
